$d = $word.ActiveDocument

# --- Step 1: Title paragraph text "1" -> "Outputs" ---
$p1 = $d.Paragraphs(1)
$p1.Range.Text = "Outputs"

# --- Step 2: Replace the old "audit page" paragraph with the new Audit section ---
$p2 = $d.Paragraphs(2)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:wordDocument xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Audit</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>Audience</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:t>Consultants</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:t>Audit team</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>Format</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t>Google Doc</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>PDF</w:t></w:r></w:p>' + `
  '</w:body></w:wordDocument>'
$p2.Range.InsertXML($xml)

# --- Step 3: Turn "Consultants" / "Audit team" into one shared bulleted list ---
$consultants = $d.Paragraphs(4)
$auditTeam = $d.Paragraphs(5)
$listRange = $d.Range($consultants.Range.Start, $auditTeam.Range.End)
$listRange.ListFormat.ApplyBulletDefault()

# --- Step 4: Bookmarks ---
# "audit" wraps the whole Audit section (Audit heading .. end of doc)
$auditHeading = $d.Paragraphs(2)
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$d.Bookmarks.Add("audit", $d.Range($auditHeading.Range.Start, $lastPara.Range.End))

# "audience" wraps Audience heading .. Audit team bullet
$audienceHeading = $d.Paragraphs(3)
$auditTeam = $d.Paragraphs(5)
$d.Bookmarks.Add("audience", $d.Range($audienceHeading.Range.Start, $auditTeam.Range.End))

# "format" wraps Format heading .. end of doc
$formatHeading = $d.Paragraphs(6)
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$d.Bookmarks.Add("format", $d.Range($formatHeading.Range.Start, $lastPara.Range.End))

Write-Output $d.Content.Text
